# Auto-generated edit script: refresh market price / profit columns (H:N)
# across all 8 sheets, matching the scheduled-runner data refresh commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 445.8
$ws.Range("I2").Value = 225
$ws.Range("J2").Value = 593
$ws.Range("K2").Value = 225
$ws.Range("L2").Value = 593
$ws.Range("M2").Value = -112
$ws.Range("N2").Value = -819
$ws.Range("H6").Value = 34.375
$ws.Range("I6").Value = 35.533333
$ws.Range("J6").Value = 17
$ws.Range("K6").Value = 106.599999
$ws.Range("L6").Value = 51
$ws.Range("M6").Value = 5.400001000000003
$ws.Range("N6").Value = -275
$ws.Range("H27").Value = 10080
$ws.Range("J27").Value = 10080
$ws.Range("L27").Value = 30240
$ws.Range("N27").Value = -30442
$ws.Range("H33").Value = 12501414
$ws.Range("I33").Value = 22727850
$ws.Range("J33").Value = 2437.111
$ws.Range("K33").Value = 22727850
$ws.Range("L33").Value = 2437.111
$ws.Range("M33").Value = -22727621
$ws.Range("N33").Value = -2895.111
$ws.Range("H40").Value = 1200
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1200
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 1200
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -1550
$ws.Range("H76").Value = 14463.917
$ws.Range("I76").Value = 8070.875
$ws.Range("K76").Value = 8070.875
$ws.Range("M76").Value = -7755.875
$ws.Range("H79").Value = 14463.917
$ws.Range("I79").Value = 8070.875
$ws.Range("K79").Value = 8070.875
$ws.Range("M79").Value = -6978.875
$ws.Range("H80").Value = 807.625
$ws.Range("I80").Value = 787.9
$ws.Range("K80").Value = 2363.7
$ws.Range("M80").Value = -1365.7
$ws.Range("H83").Value = 807.625
$ws.Range("I83").Value = 787.9
$ws.Range("K83").Value = 7091.099999999999
$ws.Range("M83").Value = -2099.099999999999
$ws.Range("H88").Value = 1013.61536
$ws.Range("J88").Value = 1034.1111
$ws.Range("L88").Value = 1034.1111
$ws.Range("N88").Value = -1846.1111
$ws.Range("H91").Value = 1013.61536
$ws.Range("J91").Value = 1034.1111
$ws.Range("L91").Value = 1034.1111
$ws.Range("N91").Value = -3842.1111
$ws.Range("H98").Value = 4980.6113
$ws.Range("I98").Value = 4978.25
$ws.Range("J98").Value = 4999.5
$ws.Range("K98").Value = 4978.25
$ws.Range("L98").Value = 4999.5
$ws.Range("M98").Value = -3480.25
$ws.Range("N98").Value = -7995.5
$ws.Range("H107").Value = 3338.5
$ws.Range("I107").Value = 2451.5
$ws.Range("J107").Value = 5999.5
$ws.Range("K107").Value = 2451.5
$ws.Range("L107").Value = 5999.5
$ws.Range("M107").Value = -531.5
$ws.Range("N107").Value = -9839.5
$ws.Range("H116").Value = 5913.1875
$ws.Range("I116").Value = 6298.6665
$ws.Range("J116").Value = 5417.5713
$ws.Range("K116").Value = 6298.6665
$ws.Range("L116").Value = 5417.5713
$ws.Range("M116").Value = -2856.6665
$ws.Range("N116").Value = -12301.5713
$ws.Range("H122").Value = 4980.6113
$ws.Range("I122").Value = 4978.25
$ws.Range("J122").Value = 4999.5
$ws.Range("K122").Value = 14934.75
$ws.Range("L122").Value = 14998.5
$ws.Range("M122").Value = -12484.75
$ws.Range("N122").Value = -19898.5
$ws.Range("H137").Value = 70855.56
$ws.Range("I137").Value = 1467.1052
$ws.Range("K137").Value = 4401.3156
$ws.Range("M137").Value = -1851.3156
$ws.Range("H138").Value = 3830.7415
$ws.Range("I138").Value = 1405.0625
$ws.Range("J138").Value = 4754.8096
$ws.Range("K138").Value = 4215.1875
$ws.Range("L138").Value = 14264.4288
$ws.Range("M138").Value = 924.8125
$ws.Range("N138").Value = -24544.4288

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1925952.9
$ws.Range("I32").Value = 2085257.2
$ws.Range("K32").Value = 2085257.2
$ws.Range("M32").Value = -2084970.2
$ws.Range("H61").Value = 1129149.4
$ws.Range("I61").Value = 3626.3044
$ws.Range("K61").Value = 3626.3044
$ws.Range("M61").Value = -3414.3044
$ws.Range("H63").Value = 4950
$ws.Range("I63").Value = 5933.3335
$ws.Range("K63").Value = 5933.3335
$ws.Range("M63").Value = -5247.3335
$ws.Range("H66").Value = 4950
$ws.Range("I66").Value = 5933.3335
$ws.Range("K66").Value = 29666.6675
$ws.Range("M66").Value = -26234.6675
$ws.Range("H74").Value = 17926.76
$ws.Range("I74").Value = 1578.0625
$ws.Range("K74").Value = 1578.0625
$ws.Range("M74").Value = -704.0625
$ws.Range("H77").Value = 17926.76
$ws.Range("I77").Value = 1578.0625
$ws.Range("K77").Value = 7890.3125
$ws.Range("M77").Value = -3522.3125
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H88").Value = 2468.9285
$ws.Range("I88").Value = 2130
$ws.Range("J88").Value = 2723.125
$ws.Range("K88").Value = 2130
$ws.Range("L88").Value = 2723.125
$ws.Range("M88").Value = -1724
$ws.Range("N88").Value = -3535.125
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H91").Value = 2468.9285
$ws.Range("I91").Value = 2130
$ws.Range("J91").Value = 2723.125
$ws.Range("K91").Value = 2130
$ws.Range("L91").Value = 2723.125
$ws.Range("M91").Value = -726
$ws.Range("N91").Value = -5531.125
$ws.Range("H102").Value = 2117.9333
$ws.Range("I102").Value = 1912.0714
$ws.Range("K102").Value = 1912.0714
$ws.Range("M102").Value = -290.0714
$ws.Range("H122").Value = 2751.6155
$ws.Range("I122").Value = 1689.1765
$ws.Range("J122").Value = 4758.4443
$ws.Range("K122").Value = 5067.529500000001
$ws.Range("L122").Value = 14275.3329
$ws.Range("M122").Value = -2617.529500000001
$ws.Range("N122").Value = -19175.3329
$ws.Range("H132").Value = 4015960.2
$ws.Range("I132").Value = 2300
$ws.Range("K132").Value = 6900
$ws.Range("M132").Value = -4370
$ws.Range("H136").Value = 1129149.4
$ws.Range("I136").Value = 3626.3044
$ws.Range("K136").Value = 10878.9132
$ws.Range("M136").Value = -8328.913199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 28878.303
$ws.Range("I20").Value = 16559.1
$ws.Range("J20").Value = 34234.477
$ws.Range("K20").Value = 16559.1
$ws.Range("L20").Value = 34234.477
$ws.Range("M20").Value = -16312.1
$ws.Range("N20").Value = -34728.477
$ws.Range("H25").Value = 1912.3334
$ws.Range("I25").Value = 1368.5
$ws.Range("K25").Value = 1368.5
$ws.Range("M25").Value = -1133.5
$ws.Range("H82").Value = 10371.667
$ws.Range("I82").Value = 10371.667
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 10371.667
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -9988.666999999999
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 10371.667
$ws.Range("I85").Value = 10371.667
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 10371.667
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -9045.666999999999
$ws.Range("N85").ClearContents()
$ws.Range("H86").Value = 2148.56
$ws.Range("I86").Value = 2135.0952
$ws.Range("K86").Value = 2135.0952
$ws.Range("M86").Value = -1012.0952
$ws.Range("H89").Value = 2148.56
$ws.Range("I89").Value = 2135.0952
$ws.Range("K89").Value = 10675.476
$ws.Range("M89").Value = -5059.476000000001
$ws.Range("H134").Value = 63882.79
$ws.Range("I134").Value = 77468.69500000001
$ws.Range("K134").Value = 232406.085
$ws.Range("M134").Value = -229871.085

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2112.4
$ws.Range("I2").Value = 1687.3334
$ws.Range("K2").Value = 1687.3334
$ws.Range("M2").Value = -1574.3334
$ws.Range("H5").Value = 1830
$ws.Range("I5").Value = 287.5
$ws.Range("K5").Value = 287.5
$ws.Range("M5").Value = -175.5
$ws.Range("H8").Value = 551.5
$ws.Range("I8").Value = 77.25
$ws.Range("J8").Value = 1500
$ws.Range("K8").Value = 77.25
$ws.Range("L8").Value = 1500
$ws.Range("M8").Value = 62.75
$ws.Range("N8").Value = -1780
$ws.Range("H58").Value = 18215.834
$ws.Range("I58").Value = 8061.2144
$ws.Range("K58").Value = 8061.2144
$ws.Range("M58").Value = -7858.2144
$ws.Range("H62").Value = 4750
$ws.Range("I62").Value = 3500
$ws.Range("K62").Value = 3500
$ws.Range("M62").Value = -2876
$ws.Range("H65").Value = 4750
$ws.Range("I65").Value = 3500
$ws.Range("K65").Value = 17500
$ws.Range("M65").Value = -14380
$ws.Range("H94").Value = 594.7857
$ws.Range("J94").Value = 643.4
$ws.Range("L94").Value = 643.4
$ws.Range("N94").Value = -1545.4
$ws.Range("H99").Value = 5775.2
$ws.Range("I99").Value = 2543.125
$ws.Range("J99").Value = 9469
$ws.Range("K99").Value = 2543.125
$ws.Range("L99").Value = 9469
$ws.Range("M99").Value = -1045.125
$ws.Range("N99").Value = -12465
$ws.Range("H126").Value = 5775.2
$ws.Range("I126").Value = 2543.125
$ws.Range("J126").Value = 9469
$ws.Range("K126").Value = 7629.375
$ws.Range("L126").Value = 28407
$ws.Range("M126").Value = -5159.375
$ws.Range("N126").Value = -33347
$ws.Range("H134").Value = 29417060
$ws.Range("I134").Value = 1757.4166
$ws.Range("J134").Value = 100013784
$ws.Range("K134").Value = 5272.2498
$ws.Range("L134").Value = 300041352
$ws.Range("M134").Value = -2737.2498
$ws.Range("N134").Value = -300046422
$ws.Range("H136").Value = 18215.834
$ws.Range("I136").Value = 8061.2144
$ws.Range("K136").Value = 24183.6432
$ws.Range("M136").Value = -21633.6432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 79134740
$ws.Range("I4").Value = 89369680
$ws.Range("J4").Value = 666831
$ws.Range("K4").Value = 268109040
$ws.Range("L4").Value = 2000493
$ws.Range("M4").Value = -268108928
$ws.Range("N4").Value = -2000717
$ws.Range("H113").Value = 57166.668
$ws.Range("I113").Value = 85000
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 255000
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = -252830
$ws.Range("N113").Value = -8840
$ws.Range("H118").Value = 13924.75
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H120").Value = 6776.0835
$ws.Range("I120").Value = 3923.6667
$ws.Range("K120").Value = 11771.0001
$ws.Range("M120").Value = -6933.000100000001
$ws.Range("H124").Value = 4464.2856
$ws.Range("I124").Value = 4208.3335
$ws.Range("K124").Value = 12625.0005
$ws.Range("M124").Value = -7715.000499999998
$ws.Range("H130").Value = 15238.429
$ws.Range("I130").Value = 2446
$ws.Range("K130").Value = 7338
$ws.Range("M130").Value = -2318
$ws.Range("H131").Value = 1445.49
$ws.Range("I131").Value = 849.8570999999999
$ws.Range("J131").Value = 1490.3226
$ws.Range("K131").Value = 2549.5713
$ws.Range("L131").Value = 4470.9678
$ws.Range("M131").Value = 2490.4287
$ws.Range("N131").Value = -14550.9678
$ws.Range("H133").Value = 37434.363
$ws.Range("I133").Value = 4129.6665
$ws.Range("K133").Value = 12388.9995
$ws.Range("M133").Value = -7328.999500000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 4473.5
$ws.Range("I46").Value = 4473.5
$ws.Range("K46").Value = 4473.5
$ws.Range("M46").Value = -4317.5
$ws.Range("H62").Value = 21000
$ws.Range("H65").Value = 21000
$ws.Range("H97").Value = 1272.7142
$ws.Range("I97").Value = 1168.8667
$ws.Range("K97").Value = 1168.8667
$ws.Range("M97").Value = -672.8667
$ws.Range("H99").Value = 4201.5557
$ws.Range("I99").Value = 4201.5557
$ws.Range("K99").Value = 4201.5557
$ws.Range("M99").Value = -1955.5557
$ws.Range("H102").Value = 8366.272000000001
$ws.Range("I102").Value = 9042.9
$ws.Range("J102").Value = 1600
$ws.Range("K102").Value = 9042.9
$ws.Range("L102").Value = 1600
$ws.Range("M102").Value = -7420.9
$ws.Range("N102").Value = -4844
$ws.Range("H140").Value = 69996
$ws.Range("J140").Value = 69996
$ws.Range("L140").Value = 69996
$ws.Range("N140").Value = -80356

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11475.056
$ws.Range("I7").Value = 17143.143
$ws.Range("J7").Value = 7868.091
$ws.Range("K7").Value = 17143.143
$ws.Range("L7").Value = 7868.091
$ws.Range("M7").Value = -17031.143
$ws.Range("N7").Value = -8092.091
$ws.Range("H40").Value = 4321
$ws.Range("I40").Value = 3648.25
$ws.Range("J40").Value = 5666.5
$ws.Range("K40").Value = 3648.25
$ws.Range("L40").Value = 5666.5
$ws.Range("M40").Value = -3512.25
$ws.Range("N40").Value = -5938.5
$ws.Range("H46").Value = 2954.7778
$ws.Range("I46").Value = 1298.75
$ws.Range("J46").Value = 3427.9285
$ws.Range("K46").Value = 1298.75
$ws.Range("L46").Value = 3427.9285
$ws.Range("M46").Value = -1110.75
$ws.Range("N46").Value = -3803.9285
$ws.Range("H63").Value = 20000
$ws.Range("J63").Value = 20333.334
$ws.Range("L63").Value = 20333.334
$ws.Range("N63").Value = -21831.334
$ws.Range("H66").Value = 20000
$ws.Range("J66").Value = 20333.334
$ws.Range("L66").Value = 61000.00199999999
$ws.Range("N66").Value = -68488.00199999999
$ws.Range("H68").Value = 12949.1
$ws.Range("I68").Value = 28124.5
$ws.Range("J68").Value = 2832.1667
$ws.Range("K68").Value = 28124.5
$ws.Range("L68").Value = 2832.1667
$ws.Range("M68").Value = -27375.5
$ws.Range("N68").Value = -4330.1667
$ws.Range("H69").Value = 176064
$ws.Range("J69").Value = 176064
$ws.Range("L69").Value = 176064
$ws.Range("N69").Value = -177686
$ws.Range("H71").Value = 12949.1
$ws.Range("I71").Value = 28124.5
$ws.Range("J71").Value = 2832.1667
$ws.Range("K71").Value = 140622.5
$ws.Range("L71").Value = 14160.8335
$ws.Range("M71").Value = -136878.5
$ws.Range("N71").Value = -21648.8335
$ws.Range("H72").Value = 176064
$ws.Range("J72").Value = 176064
$ws.Range("L72").Value = 528192
$ws.Range("N72").Value = -536304
$ws.Range("H82").Value = 2237
$ws.Range("I82").Value = 1760.4445
$ws.Range("J82").Value = 3666.6667
$ws.Range("K82").Value = 1760.4445
$ws.Range("L82").Value = 3666.6667
$ws.Range("M82").Value = -1399.4445
$ws.Range("N82").Value = -4388.6667
$ws.Range("H85").Value = 2237
$ws.Range("I85").Value = 1760.4445
$ws.Range("J85").Value = 3666.6667
$ws.Range("K85").Value = 1760.4445
$ws.Range("L85").Value = 3666.6667
$ws.Range("M85").Value = -512.4445000000001
$ws.Range("N85").Value = -6162.6667
$ws.Range("H93").Value = 7504.8887
$ws.Range("I93").Value = 18050
$ws.Range("J93").Value = 2232.3333
$ws.Range("K93").Value = 18050
$ws.Range("L93").Value = 2232.3333
$ws.Range("M93").Value = -16802
$ws.Range("N93").Value = -4728.3333
$ws.Range("H100").Value = 3122.95
$ws.Range("I100").Value = 3104.2222
$ws.Range("K100").Value = 3104.2222
$ws.Range("M100").Value = -2563.2222
$ws.Range("H122").Value = 9040.526
$ws.Range("I122").Value = 9982
$ws.Range("J122").Value = 7994.4443
$ws.Range("K122").Value = 29946
$ws.Range("L122").Value = 23983.3329
$ws.Range("M122").Value = -27496
$ws.Range("N122").Value = -28883.3329
$ws.Range("H126").Value = 11475.056
$ws.Range("I126").Value = 17143.143
$ws.Range("J126").Value = 7868.091
$ws.Range("K126").Value = 51429.429
$ws.Range("L126").Value = 23604.273
$ws.Range("M126").Value = -48959.429
$ws.Range("N126").Value = -28544.273
$ws.Range("H132").Value = 13963772
$ws.Range("J132").Value = 13963772
$ws.Range("L132").Value = 41891316
$ws.Range("N132").Value = -41896376

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H62").Value = 15972.308
$ws.Range("I62").Value = 12091.429
$ws.Range("J62").Value = 20500
$ws.Range("K62").Value = 12091.429
$ws.Range("L62").Value = 20500
$ws.Range("M62").Value = -11467.429
$ws.Range("N62").Value = -21748
$ws.Range("H65").Value = 15972.308
$ws.Range("I65").Value = 12091.429
$ws.Range("J65").Value = 20500
$ws.Range("K65").Value = 60457.145
$ws.Range("L65").Value = 102500
$ws.Range("M65").Value = -57337.145
$ws.Range("N65").Value = -108740
$ws.Range("H81").Value = 1377.1111
$ws.Range("I81").Value = 1436.0834
$ws.Range("J81").Value = 1259.1666
$ws.Range("K81").Value = 2872.1668
$ws.Range("L81").Value = 2518.3332
$ws.Range("M81").Value = -1811.1668
$ws.Range("N81").Value = -4640.3332
$ws.Range("H84").Value = 1377.1111
$ws.Range("I84").Value = 1436.0834
$ws.Range("J84").Value = 1259.1666
$ws.Range("K84").Value = 14360.834
$ws.Range("L84").Value = 12591.666
$ws.Range("M84").Value = -9056.833999999999
$ws.Range("N84").Value = -23199.666
$ws.Range("H107").Value = 1246.3158
$ws.Range("I107").Value = 1339
$ws.Range("J107").Value = 752
$ws.Range("K107").Value = 4017
$ws.Range("L107").Value = 2256
$ws.Range("M107").Value = -2097
$ws.Range("N107").Value = -6096
$ws.Range("H113").Value = 9435.714
$ws.Range("I113").Value = 10610
$ws.Range("K113").Value = 31830
$ws.Range("M113").Value = -29660
$ws.Range("H132").Value = 522676.8
$ws.Range("I132").Value = 4951.0557
$ws.Range("K132").Value = 14853.1671
$ws.Range("M132").Value = -12323.1671
$ws.Range("H136").Value = 19903.941
$ws.Range("I136").Value = 2612.6667
$ws.Range("K136").Value = 7838.000100000001
$ws.Range("M136").Value = -5288.000100000001
